# Menus, Updated Icons, Documentation
# - rename sheet to match the file name
# - simplify/update the salad menu content (ingredients, allergens, local
#   ingredients, diet codes) and add the Caesar Salad row
# - re-style the table (drop the forced centring on most data cells, keep it
#   only on the ItemName-adjacent Ingredients cell) and widen columns A/B
# - shrink the table/autofilter range down to the real data (A1:G4)
# - leave the selection sitting on C3, matching where editing wrapped up

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "salads_mtbenson"

# --- header row (text unchanged, just rewritten for safety) -----------------
$ws.Range("A1").Value = "ItemName"
$ws.Range("B1").Value = "Ingredients"
$ws.Range("C1").Value = "Allergens"
$ws.Range("D1").Value = "LocalIngredients"
$ws.Range("E1").Value = "Diet"
$ws.Range("F1").Value = "nutritionLabel"
$ws.Range("G1").Value = "LeaveEmpty"

# --- row 2: Cobb Salad --------------------------------------------------
$ws.Range("A2").Value = "Cobb Salad"
$ws.Range("B2").Value = " Egg / Guacamole / Bacon / Cheddar / Tomato / Cucumber / Romaine / Ranch Dressing"
$ws.Range("C2").Value = "Milk, eggs."
$ws.Range("D2").Value = "Fresh Start Cucumber, Fresh Start Tomato, Fresh Start Romaine"
$ws.Range("E2").Value = "BC"
$ws.Range("F2").Value = "needed"

# --- row 3: Caesar Salad (new) -------------------------------------------
$ws.Range("A3").Value = "Caesar Salad"
$ws.Range("B3").Value = "Crispy Chicken / Bacon / Romaine / Parmesan / Caesar Dressing"
$ws.Range("C3").Value = "Milk, eggs, fish."
$ws.Range("D3").Value = "Fresh Start Romaine, Castle Cheese Parmesan"
$ws.Range("E3").Value = "BC"
$ws.Range("F3").Value = "needed"

# --- row 4: House Salad ---------------------------------------------------
$ws.Range("A4").Value = "House Salad"
$ws.Range("B4").Value = "Carrots / Tomato / Cucumber / Romaine / Balsamic Dressing"
$ws.Range("C4").Value = "Dressing contains sulphites."
$ws.Range("D4").Value = "Fresh Start Cucumber, Fresh Start Tomato, Fresh Start Romaine"
$ws.Range("E4").Value = "BC, VEG, VGN, GF DF"
$ws.Range("F4").Value = "needed"

# --- styling: only the Ingredients cell for each item stays centred; every
#     other data cell reverts to the default (no alignment override) style
$ws.Range("C2:F2").Style = "Normal"
$ws.Range("B3:F3").Style = "Normal"
$ws.Range("C4:F4").Style = "Normal"

$ws.Range("B2").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B2").VerticalAlignment = -4108     # xlCenter
$ws.Range("B4").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B4").VerticalAlignment = -4108     # xlCenter

# --- column widths ---------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.8
$ws.Columns.Item(2).ColumnWidth = 73.7

# --- shrink the table/autofilter range to the real data (A1:G4) ------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G4"))

# --- selection ---------------------------------------------------------
$ws.Range("C3").Select()
